$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1040.4
$ws.Range("I6").Value = 1162.2307
$ws.Range("J6").Value = 248.5
$ws.Range("K6").Value = 3486.6921
$ws.Range("L6").Value = 745.5
$ws.Range("M6").Value = -3374.6921
$ws.Range("N6").Value = -969.5

# Row 8
$ws.Range("H8").Value = 1679
$ws.Range("I8").Value = 21
$ws.Range("J8").Value = 4995
$ws.Range("K8").Value = 63
$ws.Range("L8").Value = 14985
$ws.Range("M8").Value = 76
$ws.Range("N8").Value = -15263

# Row 38
$ws.Range("H38").Value = 840.3
$ws.Range("I38").Value = 840.3
$ws.Range("K38").Value = 2520.9
$ws.Range("M38").Value = -2148.9

# Row 39
$ws.Range("H39").Value = 1605.3846
$ws.Range("I39").Value = 22.625
$ws.Range("J39").Value = 4137.8
$ws.Range("K39").Value = 67.875
$ws.Range("L39").Value = 12413.4
$ws.Range("M39").Value = 228.125
$ws.Range("N39").Value = -13005.4

# Row 43
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 53
$ws.Range("H53").Value = 564.1818
$ws.Range("I53").Value = 530.8889
$ws.Range("J53").Value = 714
$ws.Range("K53").Value = 530.8889
$ws.Range("L53").Value = 714
$ws.Range("M53").Value = 106.1111
$ws.Range("N53").Value = -1988

# Row 55
$ws.Range("H55").Value = 234.66667
$ws.Range("I55").Value = 95.375
$ws.Range("J55").Value = 393.85715
$ws.Range("K55").Value = 95.375
$ws.Range("L55").Value = 393.85715
$ws.Range("M55").Value = 118.625
$ws.Range("N55").Value = -821.85715

# Row 103
$ws.Range("H103").Value = 1122.9286
$ws.Range("I103").Value = 1134.625
$ws.Range("J103").Value = 1107.3334
$ws.Range("K103").Value = 3403.875
$ws.Range("L103").Value = 3322.0002
$ws.Range("M103").Value = -2817.875
$ws.Range("N103").Value = -4494.0002

# Row 116
$ws.Range("H116").Value = 8087.5
$ws.Range("I116").Value = 8837.579
$ws.Range("K116").Value = 8837.579
$ws.Range("M116").Value = -5395.579

# Row 138
$ws.Range("H138").Value = 2607.6365
$ws.Range("I138").Value = 2089.111
$ws.Range("J138").Value = 2966.6155
$ws.Range("K138").Value = 6267.333
$ws.Range("L138").Value = 8899.8465
$ws.Range("M138").Value = -1127.333
$ws.Range("N138").Value = -19179.8465

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7903.449
$ws.Range("I32").Value = 4640.963
$ws.Range("J32").Value = 11907.409
$ws.Range("K32").Value = 4640.963
$ws.Range("L32").Value = 11907.409
$ws.Range("M32").Value = -4353.963
$ws.Range("N32").Value = -12481.409

# Row 61
$ws.Range("H61").Value = 125002670
$ws.Range("I61").Value = 166668960
$ws.Range("J61").Value = 3799.5
$ws.Range("K61").Value = 166668960
$ws.Range("L61").Value = 3799.5
$ws.Range("M61").Value = -166668748
$ws.Range("N61").Value = -4223.5

# Row 101
$ws.Range("H101").Value = 162017.45
$ws.Range("J101").Value = 158219.2
$ws.Range("L101").Value = 158219.2
$ws.Range("N101").Value = -164709.2

# Row 122
$ws.Range("H122").Value = 1872.2727
$ws.Range("I122").Value = 2032.0714
$ws.Range("J122").Value = 1592.625
$ws.Range("K122").Value = 6096.2142
$ws.Range("L122").Value = 4777.875
$ws.Range("M122").Value = -3646.2142
$ws.Range("N122").Value = -9677.875

# Row 136
$ws.Range("H136").Value = 125002670
$ws.Range("I136").Value = 166668960
$ws.Range("J136").Value = 3799.5
$ws.Range("K136").Value = 500006880
$ws.Range("L136").Value = 11398.5
$ws.Range("M136").Value = -500004330
$ws.Range("N136").Value = -16498.5

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2250
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2250
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2250
$ws.Range("N62").Value = -3498
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 2250
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2250
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 11250
$ws.Range("N65").Value = -17490
$ws.Range("M65").ClearContents()

# Row 99
$ws.Range("H99").Value = 3775.625
$ws.Range("I99").Value = 3934
$ws.Range("K99").Value = 3934
$ws.Range("M99").Value = -2436

# Row 107
$ws.Range("H107").Value = 661416.5600000001
$ws.Range("I107").Value = 906251
$ws.Range("K107").Value = 906251
$ws.Range("M107").Value = -904331

# Row 122
$ws.Range("H122").Value = 4123.7964
$ws.Range("I122").Value = 3996
$ws.Range("K122").Value = 11988
$ws.Range("M122").Value = -9538

# Row 126
$ws.Range("H126").Value = 3775.625
$ws.Range("I126").Value = 3934
$ws.Range("K126").Value = 11802
$ws.Range("M126").Value = -9332

# Row 134
$ws.Range("H134").Value = 50001380
$ws.Range("I134").Value = 50001380
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 150004140
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -150001605
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 108.333336
$ws.Range("I6").Value = 12.5
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 37.5
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = 75.5
$ws.Range("N6").Value = -1126

# Row 11
$ws.Range("H11").Value = 743249.5600000001
$ws.Range("I11").Value = 794460.75
$ws.Range("K11").Value = 2383382.25
$ws.Range("M11").Value = -2383242.25

# Row 59
$ws.Range("H59").Value = 6999.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 6999.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 20998.5
$ws.Range("N59").Value = -22078.5
$ws.Range("M59").ClearContents()

# Row 86
$ws.Range("H86").Value = 636.4286
$ws.Range("I86").Value = 247.5
$ws.Range("J86").Value = 677.3684
$ws.Range("K86").Value = 742.5
$ws.Range("L86").Value = 2032.1052
$ws.Range("M86").Value = 443.5
$ws.Range("N86").Value = -4404.1052

# Row 89
$ws.Range("H89").Value = 636.4286
$ws.Range("I89").Value = 247.5
$ws.Range("J89").Value = 677.3684
$ws.Range("K89").Value = 2227.5
$ws.Range("L89").Value = 6096.3156
$ws.Range("M89").Value = 3700.5
$ws.Range("N89").Value = -17952.3156

# Row 98
$ws.Range("H98").Value = 848.26666
$ws.Range("J98").Value = 576.5
$ws.Range("L98").Value = 1729.5
$ws.Range("N98").Value = -4725.5

# Row 109
$ws.Range("H109").Value = 419.57144
$ws.Range("I109").Value = 419.57144
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1258.71432
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -218.71432
$ws.Range("N109").ClearContents()

# Row 122
$ws.Range("H122").Value = 1594.2142
$ws.Range("I122").Value = 798.44446
$ws.Range("J122").Value = 3026.6
$ws.Range("K122").Value = 7186.00014
$ws.Range("L122").Value = 27239.4
$ws.Range("M122").Value = -4736.00014
$ws.Range("N122").Value = -32139.4

# Row 136
$ws.Range("H136").Value = 898.6667
$ws.Range("I136").Value = 898.6667
$ws.Range("K136").Value = 2696.0001
$ws.Range("M136").Value = 2403.9999

# Row 140
$ws.Range("H140").Value = 2127.9
$ws.Range("I140").Value = 1909.875
$ws.Range("K140").Value = 5729.625
$ws.Range("M140").Value = -549.625

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3362.375
$ws.Range("I80").Value = 3224.75
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 3224.75
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -2226.75
$ws.Range("N80").Value = -5496

# Row 83
$ws.Range("H83").Value = 3362.375
$ws.Range("I83").Value = 3224.75
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 16123.75
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -11131.75
$ws.Range("N83").Value = -27484

# Row 122
$ws.Range("H122").Value = 135831.56
$ws.Range("I122").Value = 201414.17
$ws.Range("J122").Value = 4666.3335
$ws.Range("K122").Value = 604242.51
$ws.Range("L122").Value = 13999.0005
$ws.Range("M122").Value = -601792.51
$ws.Range("N122").Value = -18899.0005

# Row 132
$ws.Range("H132").Value = 6584188.5
$ws.Range("I132").Value = 7817066.5
$ws.Range("K132").Value = 23451199.5
$ws.Range("M132").Value = -23448669.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1374.875
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 2001.3334
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 2001.3334
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -2377.3334

# Row 61
$ws.Range("H61").Value = 2288.2856
$ws.Range("I61").Value = 2009.8611
$ws.Range("K61").Value = 2009.8611
$ws.Range("M61").Value = -1807.8611

# Row 113
$ws.Range("H113").Value = 2288.2856
$ws.Range("I113").Value = 2009.8611
$ws.Range("K113").Value = 2009.8611
$ws.Range("M113").Value = 160.1388999999999

# Row 140
$ws.Range("H140").Value = 92722.5
$ws.Range("J140").Value = 92722.5
$ws.Range("L140").Value = 92722.5
$ws.Range("N140").Value = -103082.5

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3266.6667
$ws.Range("J81").Value = 2900
$ws.Range("L81").Value = 5800
$ws.Range("N81").Value = -7922

# Row 84
$ws.Range("H84").Value = 3266.6667
$ws.Range("J84").Value = 2900
$ws.Range("L84").Value = 29000
$ws.Range("N84").Value = -39608

# Row 94
$ws.Range("H94").Value = 52666.2
$ws.Range("J94").Value = 54722
$ws.Range("L94").Value = 54722
$ws.Range("N94").Value = -56524

# Row 125
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840

# Row 126
$ws.Range("H126").Value = 1698.5
$ws.Range("I126").Value = 1716.5454
$ws.Range("K126").Value = 5149.6362
$ws.Range("M126").Value = -2679.6362

# Row 136
$ws.Range("H136").Value = 62503696
$ws.Range("I136").Value = 250001700
$ws.Range("J136").Value = 4363.6665
$ws.Range("K136").Value = 750005100
$ws.Range("L136").Value = 13090.9995
$ws.Range("M136").Value = -750002550
$ws.Range("N136").Value = -18190.9995
